# custom accuracy + 데이터 1000개
# Round numeric measurement values in row 5 to 2 decimal places,
# and remove the now-superfluous row 6 (reducing the data to the
# accuracy-adjusted single row / trimmed dataset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "custom accuracy" rounding to the measurement cells on row 5.
# (Column A is the timestamp and a handful of columns already had only
# two decimals, so only the cells that actually changed are touched.)
$roundedRow5 = @{
    "C5"  = 13.45
    "D5"  = 0.82
    "F5"  = 31.77
    "G5"  = 13.66
    "H5"  = 55.42
    "I5"  = 21.29
    "J5"  = 9.44
    "L5"  = 15.51
    "N5"  = 4.47
    "Q5"  = 11.49
    "R5"  = 0.45
    "S5"  = 0.42
    "T5"  = 204.24
    "U5"  = 38.91
    "V5"  = 13.03
    "W5"  = 26.5
    "X5"  = 13.99
    "Y5"  = 1.82
    "Z5"  = 26.6
    "AA5" = 11.33
    "AE5" = 0.12
    "AF5" = 49.79
    "AG5" = 7.16
    "AH5" = 15.88
}

foreach ($addr in $roundedRow5.Keys) {
    $ws.Range($addr).Value = $roundedRow5[$addr]
}

# Drop row 6 entirely (data trimmed down), which also shrinks the sheet
# dimension from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()
